# ============================================================
# Atualizacao de bases das ligas, do dia: 27-03-2024 as 20:23
# Hungary NB I: replace row 147, and append rows 148-155
# ============================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 147 (id=145) ---
$ws.Range("A147").Value = 145
$ws.Range("B147").Value = 6818343
$ws.Range("C147").Value = 'Hungary NB I'
$ws.Range("D147").Value = 'Hungary NB I'
$ws.Range("E147").Value = 45368.35416666666
$ws.Range("F147").Value = 'Debreceni VSC'
$ws.Range("G147").Value = 'Zalaegerszegi TE'
$ws.Range("H147").Value = 5
$ws.Range("I147").Value = 1
$ws.Range("J147").Value = 'H'
$ws.Range("K147").Value = 1.727
$ws.Range("L147").Value = 3.6
$ws.Range("M147").Value = 4
$ws.Range("N147").Value = 1.666
$ws.Range("O147").Value = 3.75
$ws.Range("P147").Value = 4.333
$ws.Range("Q147").Value = -0.75
$ws.Range("R147").Value = 2
$ws.Range("S147").Value = 1.85
$ws.Range("T147").Value = 2.5
$ws.Range("U147").Value = 1.825
$ws.Range("V147").Value = 2.025
$ws.Range("W147").Value = 0.6659999999999999
$ws.Range("X147").Value = -1
$ws.Range("Y147").Value = -1
$ws.Range("Z147").Value = 1
$ws.Range("AA147").Value = -1
$ws.Range("AB147").Value = 0.825
$ws.Range("AC147").Value = -1
$ws.Range("A2").Copy()
$ws.Range("A147").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E147").PasteSpecial(-4122)

# --- Row 148 (id=146) ---
$ws.Range("A148").Value = 146
$ws.Range("B148").Value = 6818340
$ws.Range("C148").Value = 'Hungary NB I'
$ws.Range("D148").Value = 'Hungary NB I'
$ws.Range("E148").Value = 45368.44791666666
$ws.Range("F148").Value = 'Ferencvarosi TC'
$ws.Range("G148").Value = 'Puskas Academy'
$ws.Range("H148").Value = 1
$ws.Range("I148").Value = 1
$ws.Range("J148").Value = 'D'
$ws.Range("K148").Value = 1.5
$ws.Range("L148").Value = 4
$ws.Range("M148").Value = 5.5
$ws.Range("N148").Value = 1.444
$ws.Range("O148").Value = 4.2
$ws.Range("P148").Value = 6.5
$ws.Range("Q148").Value = -1.25
$ws.Range("R148").Value = 2.025
$ws.Range("S148").Value = 1.825
$ws.Range("T148").Value = 2.75
$ws.Range("U148").Value = 1.925
$ws.Range("V148").Value = 1.925
$ws.Range("W148").Value = -1
$ws.Range("X148").Value = 3.2
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = -1
$ws.Range("AA148").Value = 0.825
$ws.Range("AB148").Value = -1
$ws.Range("AC148").Value = 0.925
$ws.Range("A2").Copy()
$ws.Range("A148").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E148").PasteSpecial(-4122)

# --- Row 149 (id=147) ---
$ws.Range("A149").Value = 147
$ws.Range("B149").Value = 6818344
$ws.Range("C149").Value = 'Hungary NB I'
$ws.Range("D149").Value = 'Hungary NB I'
$ws.Range("E149").Value = 45368.5625
$ws.Range("F149").Value = 'MTK Budapest'
$ws.Range("G149").Value = 'Kisvarda FC'
$ws.Range("H149").Value = 2
$ws.Range("I149").Value = 1
$ws.Range("J149").Value = 'H'
$ws.Range("K149").Value = 1.727
$ws.Range("L149").Value = 1.2
$ws.Range("M149").Value = 4
$ws.Range("N149").Value = 2.05
$ws.Range("O149").Value = 3.4
$ws.Range("P149").Value = 3.2
$ws.Range("Q149").Value = -0.25
$ws.Range("R149").Value = 1.825
$ws.Range("S149").Value = 2.025
$ws.Range("T149").Value = 2.5
$ws.Range("U149").Value = 1.975
$ws.Range("V149").Value = 1.875
$ws.Range("W149").Value = 1.05
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = -1
$ws.Range("Z149").Value = 0.825
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = 0.9750000000000001
$ws.Range("AC149").Value = -1
$ws.Range("A2").Copy()
$ws.Range("A149").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E149").PasteSpecial(-4122)

# --- Row 150 (id=148) ---
$ws.Range("A150").Value = 148
$ws.Range("B150").Value = 6818351
$ws.Range("C150").Value = 'Hungary NB I'
$ws.Range("D150").Value = 'Hungary NB I'
$ws.Range("E150").Value = 45380.66666666666
$ws.Range("F150").Value = 'Puskas Academy'
$ws.Range("G150").Value = 'MOL Fehervar FC'
$ws.Range("K150").Value = 2.45
$ws.Range("L150").Value = 3.3
$ws.Range("M150").Value = 2.45
$ws.Range("N150").Value = 1.65
$ws.Range("O150").Value = 3.75
$ws.Range("P150").Value = 4.2
$ws.Range("Q150").Value = -0.75
$ws.Range("R150").Value = 1.925
$ws.Range("S150").Value = 1.925
$ws.Range("T150").Value = 2.5
$ws.Range("U150").Value = 1.85
$ws.Range("V150").Value = 2
$ws.Range("W150").Value = 0
$ws.Range("X150").Value = 0
$ws.Range("Y150").Value = 0
$ws.Range("Z150").Value = 0
$ws.Range("AA150").Value = 0
$ws.Range("A2").Copy()
$ws.Range("A150").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E150").PasteSpecial(-4122)

# --- Row 151 (id=149) ---
$ws.Range("A151").Value = 149
$ws.Range("B151").Value = 6818347
$ws.Range("C151").Value = 'Hungary NB I'
$ws.Range("D151").Value = 'Hungary NB I'
$ws.Range("E151").Value = 45381.4375
$ws.Range("F151").Value = 'Kisvarda FC'
$ws.Range("G151").Value = 'Debreceni VSC'
$ws.Range("K151").Value = 2.75
$ws.Range("L151").Value = 3.25
$ws.Range("M151").Value = 2.25
$ws.Range("N151").Value = 2.875
$ws.Range("O151").Value = 3.3
$ws.Range("P151").Value = 2.2
$ws.Range("Q151").Value = 0.25
$ws.Range("R151").Value = 1.825
$ws.Range("S151").Value = 2.025
$ws.Range("T151").Value = 2.25
$ws.Range("U151").Value = 1.8
$ws.Range("V151").Value = 2.05
$ws.Range("W151").Value = 0
$ws.Range("X151").Value = 0
$ws.Range("Y151").Value = 0
$ws.Range("Z151").Value = 0
$ws.Range("AA151").Value = 0
$ws.Range("A2").Copy()
$ws.Range("A151").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E151").PasteSpecial(-4122)

# --- Row 152 (id=150) ---
$ws.Range("A152").Value = 150
$ws.Range("B152").Value = 6818346
$ws.Range("C152").Value = 'Hungary NB I'
$ws.Range("D152").Value = 'Hungary NB I'
$ws.Range("E152").Value = 45381.54166666666
$ws.Range("F152").Value = 'MTK Budapest'
$ws.Range("G152").Value = 'Kecskemeti TE'
$ws.Range("K152").Value = 2.3
$ws.Range("L152").Value = 3.25
$ws.Range("M152").Value = 2.7
$ws.Range("N152").Value = 2.2
$ws.Range("O152").Value = 3.25
$ws.Range("P152").Value = 2.8
$ws.Range("Q152").Value = -0.25
$ws.Range("R152").Value = 2
$ws.Range("S152").Value = 1.85
$ws.Range("T152").Value = 2.5
$ws.Range("U152").Value = 1.9
$ws.Range("V152").Value = 1.95
$ws.Range("W152").Value = 0
$ws.Range("X152").Value = 0
$ws.Range("Y152").Value = 0
$ws.Range("Z152").Value = 0
$ws.Range("AA152").Value = 0
$ws.Range("A2").Copy()
$ws.Range("A152").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E152").PasteSpecial(-4122)

# --- Row 153 (id=151) ---
$ws.Range("A153").Value = 151
$ws.Range("B153").Value = 6818348
$ws.Range("C153").Value = 'Hungary NB I'
$ws.Range("D153").Value = 'Hungary NB I'
$ws.Range("E153").Value = 45381.64583333334
$ws.Range("F153").Value = 'Zalaegerszegi TE'
$ws.Range("G153").Value = 'Diosgyori VTK'
$ws.Range("K153").Value = 2.45
$ws.Range("L153").Value = 3.3
$ws.Range("M153").Value = 2.45
$ws.Range("N153").Value = 2.625
$ws.Range("O153").Value = 3.3
$ws.Range("P153").Value = 2.25
$ws.Range("Q153").Value = 0.25
$ws.Range("R153").Value = 1.8
$ws.Range("S153").Value = 2.05
$ws.Range("T153").Value = 2.75
$ws.Range("U153").Value = 2
$ws.Range("V153").Value = 1.85
$ws.Range("W153").Value = 0
$ws.Range("X153").Value = 0
$ws.Range("Y153").Value = 0
$ws.Range("Z153").Value = 0
$ws.Range("AA153").Value = 0
$ws.Range("A2").Copy()
$ws.Range("A153").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E153").PasteSpecial(-4122)

# --- Row 154 (id=152) ---
$ws.Range("A154").Value = 152
$ws.Range("B154").Value = 6818349
$ws.Range("C154").Value = 'Hungary NB I'
$ws.Range("D154").Value = 'Hungary NB I'
$ws.Range("E154").Value = 45382.4375
$ws.Range("F154").Value = 'Paksi'
$ws.Range("G154").Value = 'Ujpest'
$ws.Range("K154").Value = 1.666
$ws.Range("L154").Value = 3.5
$ws.Range("M154").Value = 4.333
$ws.Range("N154").Value = 1.615
$ws.Range("O154").Value = 3.5
$ws.Range("P154").Value = 4.5
$ws.Range("Q154").Value = -0.75
$ws.Range("R154").Value = 1.85
$ws.Range("S154").Value = 2
$ws.Range("T154").Value = 2.75
$ws.Range("U154").Value = 1.875
$ws.Range("V154").Value = 1.975
$ws.Range("W154").Value = 0
$ws.Range("X154").Value = 0
$ws.Range("Y154").Value = 0
$ws.Range("Z154").Value = 0
$ws.Range("AA154").Value = 0
$ws.Range("A2").Copy()
$ws.Range("A154").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E154").PasteSpecial(-4122)

# --- Row 155 (id=153) ---
$ws.Range("A155").Value = 153
$ws.Range("B155").Value = 6818350
$ws.Range("C155").Value = 'Hungary NB I'
$ws.Range("D155").Value = 'Hungary NB I'
$ws.Range("E155").Value = 45382.54166666666
$ws.Range("F155").Value = 'Mezokovesd Zsory'
$ws.Range("G155").Value = 'Ferencvarosi TC'
$ws.Range("K155").Value = 7.5
$ws.Range("L155").Value = 4.333
$ws.Range("M155").Value = 1.333
$ws.Range("N155").Value = 9.5
$ws.Range("O155").Value = 4.75
$ws.Range("P155").Value = 1.25
$ws.Range("Q155").Value = 1.5
$ws.Range("R155").Value = 1.925
$ws.Range("S155").Value = 1.925
$ws.Range("T155").Value = 2.75
$ws.Range("U155").Value = 1.825
$ws.Range("V155").Value = 2.025
$ws.Range("W155").Value = 0
$ws.Range("X155").Value = 0
$ws.Range("Y155").Value = 0
$ws.Range("Z155").Value = 0
$ws.Range("AA155").Value = 0
$ws.Range("A2").Copy()
$ws.Range("A155").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E155").PasteSpecial(-4122)

$excel.CutCopyMode = 0
"Hungary NB I sheet updated: rows 147-155 written."